$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.040112400451204
$ws.Cells.Item(2, 4).Value = 1.042756791677176
$ws.Cells.Item(2, 5).Value = 1.038298766652766
$ws.Cells.Item(2, 6).Value = 1.038860077740769
$ws.Cells.Item(2, 9).Value = 1.035744144718386
$ws.Cells.Item(2, 10).Value = 1.045200767621688
$ws.Cells.Item(2, 11).Value = 1.045532579356565
$ws.Cells.Item(2, 12).Value = 1.041087196765013
$ws.Cells.Item(2, 13).Value = 1.041646909558917
$ws.Cells.Item(2, 14).Value = 1.046685072597914
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.042139196835556
$ws.Cells.Item(3, 4).Value = 1.044698423744397
$ws.Cells.Item(3, 5).Value = 1.040074930108549
$ws.Cells.Item(3, 6).Value = 1.041528044174889
$ws.Cells.Item(3, 9).Value = 1.036269092923192
$ws.Cells.Item(3, 10).Value = 1.046867671618152
$ws.Cells.Item(3, 11).Value = 1.047282150267733
$ws.Cells.Item(3, 12).Value = 1.042670793847464
$ws.Cells.Item(3, 13).Value = 1.044120080775687
$ws.Cells.Item(3, 14).Value = 1.04835434378925
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.043443843690455
$ws.Cells.Item(4, 4).Value = 1.04594835631137
$ws.Cells.Item(4, 5).Value = 1.041217752684511
$ws.Cells.Item(4, 6).Value = 1.043246635451796
$ws.Cells.Item(4, 9).Value = 1.036604341202923
$ws.Cells.Item(4, 10).Value = 1.047939329639
$ws.Cells.Item(4, 11).Value = 1.048407420259627
$ws.Cells.Item(4, 12).Value = 1.043688593441235
$ws.Cells.Item(4, 13).Value = 1.045712411235023
$ws.Cells.Item(4, 14).Value = 1.049427523687409
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.043990718833816
$ws.Cells.Item(5, 4).Value = 1.046472321333063
$ws.Cells.Item(5, 5).Value = 1.04169667682586
$ws.Cells.Item(5, 6).Value = 1.043967324649433
$ws.Cells.Item(5, 9).Value = 1.03674423076857
$ws.Cells.Item(5, 10).Value = 1.048388223881883
$ws.Cells.Item(5, 11).Value = 1.048878882803012
$ws.Cells.Item(5, 12).Value = 1.0441148557749
$ws.Cells.Item(5, 13).Value = 1.046379965654163
$ws.Cells.Item(5, 14).Value = 1.049877055411605
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.044082448877178
$ws.Cells.Item(6, 4).Value = 1.046560209982315
$ws.Cells.Item(6, 5).Value = 1.041777002158716
$ws.Cells.Item(6, 6).Value = 1.044088227084168
$ws.Cells.Item(6, 9).Value = 1.036767657665983
$ws.Cells.Item(6, 10).Value = 1.048463500456658
$ws.Cells.Item(6, 11).Value = 1.048957950469064
$ws.Cells.Item(6, 12).Value = 1.044186332920156
$ws.Cells.Item(6, 13).Value = 1.046491943224301
$ws.Cells.Item(6, 14).Value = 1.049952438887751
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.04345115729565
$ws.Cells.Item(7, 4).Value = 1.045955363432821
$ws.Cells.Item(7, 5).Value = 1.041224158016826
$ws.Cells.Item(7, 6).Value = 1.043256272366006
$ws.Cells.Item(7, 9).Value = 1.036606214520043
$ws.Cells.Item(7, 10).Value = 1.047945334152584
$ws.Cells.Item(7, 11).Value = 1.048413726213576
$ws.Cells.Item(7, 12).Value = 1.043694295507515
$ws.Cells.Item(7, 13).Value = 1.045721338379013
$ws.Cells.Item(7, 14).Value = 1.049433536728091
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.04079880081964
$ws.Cells.Item(8, 4).Value = 1.043414327684659
$ws.Cells.Item(8, 5).Value = 1.038900388933691
$ws.Cells.Item(8, 6).Value = 1.039763367796689
$ws.Cells.Item(8, 9).Value = 1.035922477216266
$ws.Cells.Item(8, 10).Value = 1.0457655613232
$ws.Cells.Item(8, 11).Value = 1.046125286214149
$ws.Cells.Item(8, 12).Value = 1.04162382621558
$ws.Cells.Item(8, 13).Value = 1.042484410300535
$ws.Cells.Item(8, 14).Value = 1.047250668371265
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.036071107083995
$ws.Cells.Item(9, 4).Value = 1.038885905261862
$ws.Cells.Item(9, 5).Value = 1.034754601430111
$ws.Cells.Item(9, 6).Value = 1.03354663941961
$ws.Cells.Item(9, 9).Value = 1.034683212850377
$ws.Cells.Item(9, 10).Value = 1.041869971239623
$ws.Cells.Item(9, 11).Value = 1.042039097496148
$ws.Cells.Item(9, 12).Value = 1.037921247217443
$ws.Cells.Item(9, 13).Value = 1.036717241200082
$ws.Cells.Item(9, 14).Value = 1.04334954610294
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.032880772880517
$ws.Cells.Item(10, 4).Value = 1.035830653760057
$ws.Cells.Item(10, 5).Value = 1.031954429523377
$ws.Cells.Item(10, 6).Value = 1.029357302956077
$ws.Cells.Item(10, 9).Value = 1.033833154128627
$ws.Cells.Item(10, 10).Value = 1.039234256499885
$ws.Cells.Item(10, 11).Value = 1.039276857497733
$ws.Cells.Item(10, 12).Value = 1.035414544989935
$ws.Cells.Item(10, 13).Value = 1.032826804209674
$ws.Cells.Item(10, 14).Value = 1.040710088345951
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.031489663312196
$ws.Cells.Item(11, 4).Value = 1.034498599105757
$ws.Cells.Item(11, 5).Value = 1.030732855434125
$ws.Cells.Item(11, 6).Value = 1.027531888360049
$ws.Cells.Item(11, 9).Value = 1.033459238631617
$ws.Cells.Item(11, 10).Value = 1.038083347720465
$ws.Cells.Item(11, 11).Value = 1.038071277542722
$ws.Cells.Item(11, 12).Value = 1.034319595633216
$ws.Cells.Item(11, 13).Value = 1.031130671138465
$ws.Cells.Item(11, 14).Value = 1.039557545144052
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.030971444720128
$ws.Cells.Item(12, 4).Value = 1.034002403571806
$ws.Cells.Item(12, 5).Value = 1.03027770440337
$ws.Cells.Item(12, 6).Value = 1.02685206968833
$ws.Cells.Item(12, 9).Value = 1.033319458668881
$ws.Cells.Item(12, 10).Value = 1.037654363681859
$ws.Cells.Item(12, 11).Value = 1.037622002401601
$ws.Cells.Item(12, 12).Value = 1.033911413209556
$ws.Cells.Item(12, 13).Value = 1.030498855754735
$ws.Cells.Item(12, 14).Value = 1.039127951898906
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.031082672915135
$ws.Cells.Item(13, 4).Value = 1.034108903743732
$ws.Cells.Item(13, 5).Value = 1.030375400024907
$ws.Cells.Item(13, 6).Value = 1.02699797463663
$ws.Cells.Item(13, 9).Value = 1.033349482506117
$ws.Cells.Item(13, 10).Value = 1.037746450059624
$ws.Cells.Item(13, 11).Value = 1.037718440574009
$ws.Cells.Item(13, 12).Value = 1.033999036833555
$ws.Cells.Item(13, 13).Value = 1.03063446456741
$ws.Cells.Item(13, 14).Value = 1.03922016904989
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.031446858036551
$ws.Cells.Item(14, 4).Value = 1.034457612462851
$ws.Cells.Item(14, 5).Value = 1.030695261382021
$ws.Cells.Item(14, 6).Value = 1.027475730996929
$ws.Cells.Item(14, 9).Value = 1.033447702639666
$ws.Cells.Item(14, 10).Value = 1.038047918293307
$ws.Cells.Item(14, 11).Value = 1.038034170522448
$ws.Cells.Item(14, 12).Value = 1.034285885343907
$ws.Cells.Item(14, 13).Value = 1.031078482078262
$ws.Cells.Item(14, 14).Value = 1.039522065403044
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.031671044568452
$ws.Cells.Item(15, 4).Value = 1.03467227514354
$ws.Cells.Item(15, 5).Value = 1.03089215127337
$ws.Cells.Item(15, 6).Value = 1.027769854831814
$ws.Cells.Item(15, 9).Value = 1.033508100786118
$ws.Cells.Item(15, 10).Value = 1.038233464829898
$ws.Cells.Item(15, 11).Value = 1.038228506286869
$ws.Cells.Item(15, 12).Value = 1.03446242632028
$ws.Cells.Item(15, 13).Value = 1.031351816079777
$ws.Cells.Item(15, 14).Value = 1.039707875437
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.03297288858951
$ws.Cells.Item(16, 4).Value = 1.035918862294379
$ws.Cells.Item(16, 5).Value = 1.032035306673764
$ws.Cells.Item(16, 6).Value = 1.029478203914281
$ws.Cells.Item(16, 9).Value = 1.033857845509616
$ws.Cells.Item(16, 10).Value = 1.03931043239163
$ws.Cells.Item(16, 11).Value = 1.039356664068093
$ws.Cells.Item(16, 12).Value = 1.035487009179078
$ws.Cells.Item(16, 13).Value = 1.032939122485777
$ws.Cells.Item(16, 14).Value = 1.040786372416199
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.033786879388087
$ws.Cells.Item(17, 4).Value = 1.036698344678012
$ws.Cells.Item(17, 5).Value = 1.032749918802562
$ws.Cells.Item(17, 6).Value = 1.030546707412703
$ws.Cells.Item(17, 9).Value = 1.034075658919049
$ws.Cells.Item(17, 10).Value = 1.039983381052987
$ws.Cells.Item(17, 11).Value = 1.040061753056218
$ws.Cells.Item(17, 12).Value = 1.036127124982388
$ws.Cells.Item(17, 13).Value = 1.033931663374258
$ws.Cells.Item(17, 14).Value = 1.041460276741844
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.034260736164394
$ws.Cells.Item(18, 4).Value = 1.03715212742867
$ws.Cells.Item(18, 5).Value = 1.033165866484954
$ws.Cells.Item(18, 6).Value = 1.031168851540404
$ws.Cells.Item(18, 9).Value = 1.034202143852595
$ws.Cells.Item(18, 10).Value = 1.040374974435058
$ws.Cells.Item(18, 11).Value = 1.040472104526565
$ws.Cells.Item(18, 12).Value = 1.036499576698134
$ws.Cells.Item(18, 13).Value = 1.034509485636661
$ws.Cells.Item(18, 14).Value = 1.041852426231433
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.034422152640892
$ws.Cells.Item(19, 4).Value = 1.0373067082813
$ws.Cells.Item(19, 5).Value = 1.033307546926975
$ws.Cells.Item(19, 6).Value = 1.031380802678088
$ws.Cells.Item(19, 9).Value = 1.034245177066175
$ws.Cells.Item(19, 10).Value = 1.040508341727986
$ws.Cells.Item(19, 11).Value = 1.040611869771632
$ws.Cells.Item(19, 12).Value = 1.036626418695158
$ws.Cells.Item(19, 13).Value = 1.034706321642136
$ws.Cells.Item(19, 14).Value = 1.041985982921215
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.033699642426883
$ws.Cells.Item(20, 4).Value = 1.036614804522053
$ws.Cells.Item(20, 5).Value = 1.032673338317299
$ws.Cells.Item(20, 6).Value = 1.030432180918712
$ws.Cells.Item(20, 9).Value = 1.034052347830276
$ws.Cells.Item(20, 10).Value = 1.039911276109612
$ws.Cells.Item(20, 11).Value = 1.039986198611408
$ws.Cells.Item(20, 12).Value = 1.036058541734345
$ws.Cells.Item(20, 13).Value = 1.033825288385411
$ws.Cells.Item(20, 14).Value = 1.041388069401178
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.031339656330277
$ws.Cells.Item(21, 4).Value = 1.034354965734823
$ws.Cells.Item(21, 5).Value = 1.030601109265189
$ws.Cells.Item(21, 6).Value = 1.027335093278613
$ws.Cells.Item(21, 9).Value = 1.033418803975374
$ws.Cells.Item(21, 10).Value = 1.037959184729858
$ws.Cells.Item(21, 11).Value = 1.037941236800912
$ws.Cells.Item(21, 12).Value = 1.034201456477589
$ws.Cells.Item(21, 13).Value = 1.030947780064264
$ws.Cells.Item(21, 14).Value = 1.039433205827756
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.029847146154464
$ws.Cells.Item(22, 4).Value = 1.032925929446692
$ws.Cells.Item(22, 5).Value = 1.02929007313055
$ws.Cells.Item(22, 6).Value = 1.025377508787243
$ws.Cells.Item(22, 9).Value = 1.033015306757719
$ws.Cells.Item(22, 10).Value = 1.036723214194466
$ws.Cells.Item(22, 11).Value = 1.036646968151986
$ws.Cells.Item(22, 12).Value = 1.03302531213702
$ws.Cells.Item(22, 13).Value = 1.029128153864921
$ws.Cells.Item(22, 14).Value = 1.038195480072436
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.030639193381334
$ws.Cells.Item(23, 4).Value = 1.03368427892366
$ws.Cells.Item(23, 5).Value = 1.029985863555686
$ws.Cells.Item(23, 6).Value = 1.026416261899555
$ws.Cells.Item(23, 9).Value = 1.033229702600125
$ws.Cells.Item(23, 10).Value = 1.037379255202918
$ws.Cells.Item(23, 11).Value = 1.037333905625778
$ws.Cells.Item(23, 12).Value = 1.033649628959656
$ws.Cells.Item(23, 13).Value = 1.030093780793315
$ws.Cells.Item(23, 14).Value = 1.03885245273437
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.033739063927131
$ws.Cells.Item(24, 4).Value = 1.036652555427969
$ws.Cells.Item(24, 5).Value = 1.032707944439513
$ws.Cells.Item(24, 6).Value = 1.030483933894757
$ws.Cells.Item(24, 9).Value = 1.034062882843638
$ws.Cells.Item(24, 10).Value = 1.039943860089964
$ws.Cells.Item(24, 11).Value = 1.040020341238391
$ws.Cells.Item(24, 12).Value = 1.036089534383839
$ws.Cells.Item(24, 13).Value = 1.033873358091113
$ws.Cells.Item(24, 14).Value = 1.041420699654519
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.037299950792748
$ws.Cells.Item(25, 4).Value = 1.040062849326943
$ws.Cells.Item(25, 5).Value = 1.035832637134467
$ws.Cells.Item(25, 6).Value = 1.035161471349466
$ws.Cells.Item(25, 9).Value = 1.035007747212139
$ws.Cells.Item(25, 10).Value = 1.042883740514945
$ws.Cells.Item(25, 11).Value = 1.043102041676549
$ws.Cells.Item(25, 12).Value = 1.03888506385237
$ws.Cells.Item(25, 13).Value = 1.038216008700185
$ws.Cells.Item(25, 14).Value = 1.044364755046913
